# Refresh market-board derived columns (H:N) on the per-job Leve Profit
# sheets, as produced by the scheduled FFXIV market-data pull.
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18: "You Grow, Girl" (Leve Item ID 5471)
$ws.Cells.Item(18, 8).Value = 1242.7778  # H18
$ws.Cells.Item(18, 9).Value = 270.4  # I18
$ws.Cells.Item(18, 10).Value = 2458.25  # J18
$ws.Cells.Item(18, 11).Value = 270.4  # K18
$ws.Cells.Item(18, 12).Value = 2458.25  # L18
$ws.Cells.Item(18, 13).Value = 13.60000000000002  # M18
$ws.Cells.Item(18, 14).Value = -3026.25  # N18
# Row 29: "Dripping with Venom" (Leve Item ID 4575)
$ws.Cells.Item(29, 8).Value = 350  # H29
$ws.Cells.Item(29, 9).Value = 200  # I29
$ws.Cells.Item(29, 10).Value = 500  # J29
$ws.Cells.Item(29, 11).Value = 600  # K29
$ws.Cells.Item(29, 12).Value = 1500  # L29
$ws.Cells.Item(29, 13).Value = -319  # M29
$ws.Cells.Item(29, 14).Value = -2062  # N29
# Row 38: "Just Give Him a Serum" (Leve Item ID 4599)
$ws.Cells.Item(38, 8).Value = 1457.1111  # H38
$ws.Cells.Item(38, 9).Value = 1599.4286  # I38
$ws.Cells.Item(38, 10).Value = 959  # J38
$ws.Cells.Item(38, 11).Value = 4798.2858  # K38
$ws.Cells.Item(38, 12).Value = 2877  # L38
$ws.Cells.Item(38, 13).Value = -4426.2858  # M38
$ws.Cells.Item(38, 14).Value = -3621  # N38
# Row 58: "A Matter of Vital Importance" (Leve Item ID 4606)
$ws.Cells.Item(58, 8).Value = 6719.1665  # H58
$ws.Cells.Item(58, 9).Value = 215  # I58
$ws.Cells.Item(58, 10).Value = 8020  # J58
$ws.Cells.Item(58, 11).Value = 645  # K58
$ws.Cells.Item(58, 12).Value = 24060  # L58
$ws.Cells.Item(58, 13).Value = -495  # M58
$ws.Cells.Item(58, 14).Value = -24360  # N58
# Row 76: "Warding Off Temptation" (Leve Item ID 12602)
$ws.Cells.Item(76, 8).Value = 3089666.2  # H76
$ws.Cells.Item(76, 9).Value = 4118307.8  # I76
$ws.Cells.Item(76, 10).Value = 3742.2222  # J76
$ws.Cells.Item(76, 11).Value = 4118307.8  # K76
$ws.Cells.Item(76, 12).Value = 3742.2222  # L76
$ws.Cells.Item(76, 13).Value = -4117992.8  # M76
$ws.Cells.Item(76, 14).Value = -4372.2222  # N76
# Row 79: "The Garden of Arcane Delights (L)" (Leve Item ID 12602)
$ws.Cells.Item(79, 8).Value = 3089666.2  # H79
$ws.Cells.Item(79, 9).Value = 4118307.8  # I79
$ws.Cells.Item(79, 10).Value = 3742.2222  # J79
$ws.Cells.Item(79, 11).Value = 4118307.8  # K79
$ws.Cells.Item(79, 12).Value = 3742.2222  # L79
$ws.Cells.Item(79, 13).Value = -4117215.8  # M79
$ws.Cells.Item(79, 14).Value = -5926.2222  # N79
# Row 129: "Practical Command" (Leve Item ID 36115)
$ws.Cells.Item(129, 8).Value = 1162.5625  # H129
$ws.Cells.Item(129, 9).Value = 395  # I129
$ws.Cells.Item(129, 10).Value = 1213.7333  # J129
$ws.Cells.Item(129, 11).Value = 1185  # K129
$ws.Cells.Item(129, 12).Value = 3641.199900000001  # L129
$ws.Cells.Item(129, 13).Value = 3815  # M129
$ws.Cells.Item(129, 14).Value = -13641.1999  # N129
# Row 134: "Binding Spells" (Leve Item ID 41997)
$ws.Cells.Item(134, 8).Value = 77780  # H134
$ws.Cells.Item(134, 9).Value = 0  # I134
$ws.Cells.Item(134, 10).Value = 77780  # J134
$ws.Cells.Item(134, 11).Value = 0  # K134
$ws.Cells.Item(134, 12).Value = 77780  # L134
$ws.Cells.Item(134, 14).Value = -87920  # N134
# Row 135: "For Tired Minds" (Leve Item ID 44047)
$ws.Cells.Item(135, 8).Value = 2367.6667  # H135
$ws.Cells.Item(135, 9).Value = 2367.6667  # I135
$ws.Cells.Item(135, 10).Value = 0  # J135
$ws.Cells.Item(135, 11).Value = 21309.0003  # K135
$ws.Cells.Item(135, 12).Value = 0  # L135
$ws.Cells.Item(135, 13).Value = -18774.0003  # M135
$ws.Cells.Item(135, 14).ClearContents()  # N135

$ws = $wb.Worksheets.Item("ARM")
# Row 2: "Ain't Got No Ingots" (Leve Item ID 27713)
$ws.Cells.Item(2, 8).Value = 35333.38  # H2
$ws.Cells.Item(2, 9).Value = 37921.035  # I2
$ws.Cells.Item(2, 10).Value = 400  # J2
$ws.Cells.Item(2, 11).Value = 37921.035  # K2
$ws.Cells.Item(2, 12).Value = 400  # L2
$ws.Cells.Item(2, 13).Value = -37808.035  # M2
$ws.Cells.Item(2, 14).Value = -626  # N2
# Row 32: "Ingot We Trust" (Leve Item ID 44147)
$ws.Cells.Item(32, 8).Value = 5586.35  # H32
$ws.Cells.Item(32, 9).Value = 3139.6162  # I32
$ws.Cells.Item(32, 10).Value = 20616.285  # J32
$ws.Cells.Item(32, 11).Value = 3139.6162  # K32
$ws.Cells.Item(32, 12).Value = 20616.285  # L32
$ws.Cells.Item(32, 13).Value = -2852.6162  # M32
$ws.Cells.Item(32, 14).Value = -21190.285  # N32
# Row 97: "Ore for Me" (Leve Item ID 19941)
$ws.Cells.Item(97, 8).Value = 5647.8423  # H97
$ws.Cells.Item(97, 9).Value = 6258.706  # I97
$ws.Cells.Item(97, 10).Value = 455.5  # J97
$ws.Cells.Item(97, 11).Value = 6258.706  # K97
$ws.Cells.Item(97, 12).Value = 455.5  # L97
$ws.Cells.Item(97, 13).Value = -5762.706  # M97
$ws.Cells.Item(97, 14).Value = -1447.5  # N97
# Row 116: "No Scope" (Leve Item ID 27713)
$ws.Cells.Item(116, 8).Value = 35333.38  # H116
$ws.Cells.Item(116, 9).Value = 37921.035  # I116
$ws.Cells.Item(116, 10).Value = 400  # J116
$ws.Cells.Item(116, 11).Value = 37921.035  # K116
$ws.Cells.Item(116, 12).Value = 400  # L116
$ws.Cells.Item(116, 13).Value = -35627.035  # M116
$ws.Cells.Item(116, 14).Value = -4988  # N116

$ws = $wb.Worksheets.Item("BSM")
# Row 3: "Hells Bells" (Leve Item ID 27713)
$ws.Cells.Item(3, 8).Value = 35333.38  # H3
$ws.Cells.Item(3, 9).Value = 37921.035  # I3
$ws.Cells.Item(3, 10).Value = 400  # J3
$ws.Cells.Item(3, 11).Value = 37921.035  # K3
$ws.Cells.Item(3, 12).Value = 400  # L3
$ws.Cells.Item(3, 13).Value = -37807.035  # M3
$ws.Cells.Item(3, 14).Value = -628  # N3
# Row 54: "Get Me to the War on Time" (Leve Item ID 2376)
$ws.Cells.Item(54, 8).Value = 3150  # H54
$ws.Cells.Item(54, 9).Value = 3150  # I54
$ws.Cells.Item(54, 10).Value = 0  # J54
$ws.Cells.Item(54, 11).Value = 3150  # K54
$ws.Cells.Item(54, 12).Value = 0  # L54
$ws.Cells.Item(54, 13).Value = -2666  # M54
# Row 86: "Through Thick and Thin" (Leve Item ID 12526)
$ws.Cells.Item(86, 8).Value = 1435.6  # H86
$ws.Cells.Item(86, 9).Value = 1277.5555  # I86
$ws.Cells.Item(86, 10).Value = 1564.909  # J86
$ws.Cells.Item(86, 11).Value = 1277.5555  # K86
$ws.Cells.Item(86, 12).Value = 1564.909  # L86
$ws.Cells.Item(86, 13).Value = -154.5554999999999  # M86
$ws.Cells.Item(86, 14).Value = -3810.909  # N86
# Row 89: "Piercing Eyes Deserve Piercing Shafts (L)" (Leve Item ID 12526)
$ws.Cells.Item(89, 8).Value = 1435.6  # H89
$ws.Cells.Item(89, 9).Value = 1277.5555  # I89
$ws.Cells.Item(89, 10).Value = 1564.909  # J89
$ws.Cells.Item(89, 11).Value = 6387.7775  # K89
$ws.Cells.Item(89, 12).Value = 7824.545  # L89
$ws.Cells.Item(89, 13).Value = -771.7775000000001  # M89
$ws.Cells.Item(89, 14).Value = -19056.545  # N89
# Row 94: "High Steal" (Leve Item ID 19939)
$ws.Cells.Item(94, 8).Value = 1000.56665  # H94
$ws.Cells.Item(94, 9).Value = 603.88464  # I94
$ws.Cells.Item(94, 10).Value = 3579  # J94
$ws.Cells.Item(94, 11).Value = 603.88464  # K94
$ws.Cells.Item(94, 12).Value = 3579  # L94
$ws.Cells.Item(94, 13).Value = -152.88464  # M94
$ws.Cells.Item(94, 14).Value = -4481  # N94
# Row 105: "Ingot to Wing It" (Leve Item ID 19947)
$ws.Cells.Item(105, 8).Value = 247002.97  # H105
$ws.Cells.Item(105, 9).Value = 2866.6667  # I105
$ws.Cells.Item(105, 10).Value = 717837.3  # J105
$ws.Cells.Item(105, 11).Value = 2866.6667  # K105
$ws.Cells.Item(105, 12).Value = 717837.3  # L105
$ws.Cells.Item(105, 13).Value = -1119.6667  # M105
$ws.Cells.Item(105, 14).Value = -721331.3  # N105

$ws = $wb.Worksheets.Item("CRP")
# Row 54: "The Turning Point" (Leve Item ID 2413)
$ws.Cells.Item(54, 8).Value = 4872.5  # H54
$ws.Cells.Item(54, 9).Value = 2345  # I54
$ws.Cells.Item(54, 10).Value = 7400  # J54
$ws.Cells.Item(54, 11).Value = 2345  # K54
$ws.Cells.Item(54, 12).Value = 7400  # L54
$ws.Cells.Item(54, 13).Value = -1687  # M54
$ws.Cells.Item(54, 14).Value = -8716  # N54
# Row 94: "Beech, Please" (Leve Item ID 32934)
$ws.Cells.Item(94, 8).Value = 2133.5715  # H94
$ws.Cells.Item(94, 9).Value = 1100  # I94
$ws.Cells.Item(94, 10).Value = 2547  # J94
$ws.Cells.Item(94, 11).Value = 1100  # K94
$ws.Cells.Item(94, 12).Value = 2547  # L94
$ws.Cells.Item(94, 13).Value = -649  # M94
$ws.Cells.Item(94, 14).Value = -3449  # N94
# Row 134: "Wood You Be Quiet" (Leve Item ID 44020)
$ws.Cells.Item(134, 8).Value = 2558.1177  # H134
$ws.Cells.Item(134, 9).Value = 1034.76  # I134
$ws.Cells.Item(134, 10).Value = 6789.6665  # J134
$ws.Cells.Item(134, 11).Value = 3104.28  # K134
$ws.Cells.Item(134, 12).Value = 20368.9995  # L134
$ws.Cells.Item(134, 13).Value = -569.2799999999997  # M134
$ws.Cells.Item(134, 14).Value = -25438.9995  # N134
# Row 135: "The Wing's Wings" (Leve Item ID 42008)
$ws.Cells.Item(135, 8).Value = 43842.5  # H135
$ws.Cells.Item(135, 9).Value = 0  # I135
$ws.Cells.Item(135, 10).Value = 43842.5  # J135
$ws.Cells.Item(135, 11).Value = 0  # K135
$ws.Cells.Item(135, 12).Value = 43842.5  # L135
$ws.Cells.Item(135, 14).Value = -53982.5  # N135

$ws = $wb.Worksheets.Item("CUL")
# Row 5: "What a Sap" (Leve Item ID 43974)
$ws.Cells.Item(5, 8).Value = 1444.4103  # H5
$ws.Cells.Item(5, 9).Value = 834.6  # I5
$ws.Cells.Item(5, 10).Value = 1654.6897  # J5
$ws.Cells.Item(5, 11).Value = 2503.8  # K5
$ws.Cells.Item(5, 12).Value = 4964.0691  # L5
$ws.Cells.Item(5, 13).Value = -2391.8  # M5
$ws.Cells.Item(5, 14).Value = -5188.0691  # N5
# Row 131: "The Mountain Steeped" (Leve Item ID 36060)
$ws.Cells.Item(131, 8).Value = 4903343.5  # H131
$ws.Cells.Item(131, 9).Value = 516.6667  # I131
$ws.Cells.Item(131, 10).Value = 5129628  # J131
$ws.Cells.Item(131, 11).Value = 1550.0001  # K131
$ws.Cells.Item(131, 12).Value = 15388884  # L131
$ws.Cells.Item(131, 13).Value = 3489.9999  # M131
$ws.Cells.Item(131, 14).Value = -15398964  # N131
# Row 135: "Not-so-secret Ingredient" (Leve Item ID 43974)
$ws.Cells.Item(135, 8).Value = 1444.4103  # H135
$ws.Cells.Item(135, 9).Value = 834.6  # I135
$ws.Cells.Item(135, 10).Value = 1654.6897  # J135
$ws.Cells.Item(135, 11).Value = 7511.400000000001  # K135
$ws.Cells.Item(135, 12).Value = 14892.2073  # L135
$ws.Cells.Item(135, 13).Value = -4976.400000000001  # M135
$ws.Cells.Item(135, 14).Value = -19962.2073  # N135

$ws = $wb.Worksheets.Item("LTW")
# Row 16: "Saddle Sore" (Leve Item ID 5289)
$ws.Cells.Item(16, 8).Value = 7694047  # H16
$ws.Cells.Item(16, 9).Value = 11112091  # I16
$ws.Cells.Item(16, 10).Value = 3448  # J16
$ws.Cells.Item(16, 11).Value = 11112091  # K16
$ws.Cells.Item(16, 12).Value = 3448  # L16
$ws.Cells.Item(16, 13).Value = -11111921  # M16
$ws.Cells.Item(16, 14).Value = -3788  # N16
# Row 55: "It's Not a Job, It's a Calling" (Leve Item ID 5284)
$ws.Cells.Item(55, 8).Value = 414.27777  # H55
$ws.Cells.Item(55, 9).Value = 305.3  # I55
$ws.Cells.Item(55, 10).Value = 550.5  # J55
$ws.Cells.Item(55, 11).Value = 305.3  # K55
$ws.Cells.Item(55, 12).Value = 550.5  # L55
$ws.Cells.Item(55, 13).Value = -132.3  # M55
$ws.Cells.Item(55, 14).Value = -896.5  # N55
